$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 2 (latest date goes on top); this shifts
# the previous rows 2-5 down to rows 3-6.
$ws.Rows.Item(2).Insert()

# The Insert() call copies formatting from the row above (the bold header),
# so strip that back off to match the unstyled data rows.
$ws.Range("A2:D2").Style = "Normal"

# Write the new date as literal text (matching the other date cells, which
# are stored as text, not as real Excel dates) rather than letting Excel
# auto-convert the string into a date serial number.
$dateCell = $ws.Cells.Item(2, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-11-25"
$dateCell.Style = "Normal"

# The new row carries the same commodity prices as the rest of the table.
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
